$d = $word.ActiveDocument

# Locate the "save playlists from other people" paragraph and insert
# three new ListParagraph (numId=1, ilvl=0) bullet items right after it,
# before the "literally unrelated..." heading.
$r = $d.Content
[void]$r.Find.Execute("save playlists from other people", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$anchor = $r.Paragraphs(1)

$newLines = @(
    "routing to view a current playlist",
    "save listened videos to a mongodb db",
    "volume controls app wide"
)

foreach ($line in $newLines) {
    $anchor.Range.InsertParagraphAfter()
    $anchor = $anchor.Next()
    $anchor.Range.Text = $line
}
